$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New gift-certificate codes for rows 2-11
$codes = @(
    "PGAGHPSCYCFPDRDP ",
    "RPALZCIOTDMTCJTK ",
    "GLOIWCSCSCTVGYIF ",
    "ISWGFHJZJRDMTCIL ",
    "VSQGZPMZZYLGVIVF ",
    "HIZYMCWGPLDQLYJZ",
    "PTJAMJQLAGAYMGZS",
    "SAHISDQPYFFHGKTA ",
    "FPHKLTYHVAYTHTHC ",
    "MOFLQCWGOPMSJFHV "
)

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $codes[$i]
}

# Re-apply the per-row look-and-feel that the new data uses (do this before the
# trailing rows are removed, since some source styles live on those rows).

# Row 6 takes on the "bold / dark grey" emphasis style that used to live on row 13.
$ws.Range("A13").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Row 8 goes back to the regular body style (same as row 7/others).
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Rows 9-11 take on the "bold / plain black, unprotected" emphasis style that
# used to live on rows 3 / 11 / 12 / 14.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A9:A11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

# Remove the now-unused trailing rows (12-18) so dimension becomes A1:A11
$ws.Range("A12:A18").EntireRow.Delete()

$ws.Range("A11").Select()
